$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: remember the original ReferenciaMB (Q) / EntidadeMB (R) values,
# which are tied to the row POSITION and must stay exactly where they are,
# even though three data rows below them are going to be removed.
$qrOriginal = @{}
for ($r = 2; $r -le 40; $r++) {
    $q = $ws.Cells.Item($r, 17).Value2
    $rr = $ws.Cells.Item($r, 18).Value2
    $qrOriginal[$r] = @($q, $rr)
}

# --- Step 2: remove the three rows whose orders are no longer present
# (NumeroPedido 748322, 748361, 748679). Delete bottom-up so the row
# numbers of the earlier deletions stay valid.
$ws.Rows.Item(39).Delete()
$ws.Rows.Item(24).Delete()
$ws.Rows.Item(19).Delete()

# --- Step 3: put the ReferenciaMB / EntidadeMB values back where they
# belong (tied to the row position, unaffected by the deletions above).
for ($r = 2; $r -le 37; $r++) {
    $vals = $qrOriginal[$r]
    $cQ = $ws.Cells.Item($r, 17)
    $cQ.NumberFormat = "@"
    $cQ.Value = $vals[0]
    $cQ.Style = "Normal"

    $cR = $ws.Cells.Item($r, 18)
    $cR.NumberFormat = "@"
    $cR.Value = $vals[1]
    $cR.Style = "Normal"
}

# --- Step 4: fill in the now-computed IVA / ValorTotal / MontanteMB
# amounts (previously placeholder zeros) for every remaining data row.
for ($r = 2; $r -le 37; $r++) {
    $cN = $ws.Cells.Item($r, 14)
    $cN.NumberFormat = "@"
    $cN.Value = "28€"
    $cN.Style = "Normal"

    $cO = $ws.Cells.Item($r, 15)
    $cO.NumberFormat = "@"
    $cO.Value = "6,44€"
    $cO.Style = "Normal"

    $cP = $ws.Cells.Item($r, 16)
    $cP.NumberFormat = "@"
    $cP.Value = "34,44€"
    $cP.Style = "Normal"
}
